$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.407.30'
$ws.Range("E2").Value = '  -2.26%  '
$ws.Range("D3").Value = '2.638.68'
$ws.Range("E3").Value = '  -3.52%  '
$ws.Range("E4").Value = '  +0.08%  '
$__style = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '598.09'
$ws.Range("D5").Style = $__style
$ws.Range("E5").Value = '  -0.99%  '
$__style = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '166.86'
$ws.Range("D6").Style = $__style
$ws.Range("E6").Value = '  -1.73%  '
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("E8").Value = '  -0.81%  '
$ws.Range("D9").Value = '2.638.09'
$ws.Range("E9").Value = '  -3.43%  '
$__style = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.144'
$ws.Range("D10").Style = $__style
$ws.Range("E10").Value = '  -1.41%  '
$ws.Range("E11").Value = '  +1.47%  '
$ws.Range("E12").Value = '  -1.63%  '
$ws.Range("E13").Value = '  -2.23%  '
$__style = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.92'
$ws.Range("D14").Style = $__style
$ws.Range("E14").Value = '  -2.86%  '
$ws.Range("D15").Value = '3.120.83'
$ws.Range("E15").Value = '  -3.44%  '
$ws.Range("E16").Value = '  -3.71%  '
$ws.Range("D17").Value = '67.297.20'
$ws.Range("E17").Value = '  -2.16%  '
$ws.Range("D18").Value = '2.632.08'
$ws.Range("E18").Value = '  -3.99%  '
$__style = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.84'
$ws.Range("D19").Style = $__style
$ws.Range("E19").Value = '  -1.09%  '
$__style = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.84'
$ws.Range("D20").Style = $__style
$ws.Range("E20").Value = '  +2.04%  '
$__style = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '362.74'
$ws.Range("D21").Style = $__style
$ws.Range("E21").Value = '  -3.04%  '
$__style = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.40'
$ws.Range("D22").Style = $__style
$ws.Range("E22").Value = '  -3.24%  '
$ws.Range("E23").Value = '  -4.25%  '
$__style = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '10.95'
$ws.Range("D24").Style = $__style
$ws.Range("E24").Value = '  +8.16%  '
$ws.Range("E25").Value = '  -6.10%  '
$ws.Range("E26").Value = '  +0.17%  '
$__style = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '70.81'
$ws.Range("D27").Style = $__style
$ws.Range("E27").Value = '  -4.07%  '
$ws.Range("E29").Value = '  -3.81%  '
$ws.Range("E30").Value = '  +0.07%  '
$__style = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '553.47'
$ws.Range("D31").Style = $__style
$ws.Range("E31").Value = '  -5.63%  '
$__style = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.04'
$ws.Range("D32").Style = $__style
$ws.Range("E32").Value = '  -3.61%  '
$ws.Range("E33").Value = '  -4.45%  '
$ws.Range("E34").Value = '  -2.09%  '
$ws.Range("E35").Value = '  +0.06%  '
$ws.Range("E36").Value = '  +0.03%  '
$ws.Range("E37").Value = '  -5.70%  '
$ws.Range("E38").Value = '  -2.77%  '
$__style = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '19.38'
$ws.Range("D39").Style = $__style
$ws.Range("E39").Value = '  -3.26%  '
$ws.Range("E40").Value = '  -2.74%  '
$ws.Range("E41").Value = '  -5.31%  '
$ws.Range("E42").Value = '  -4.54%  '
$__style = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '17.92'
$ws.Range("D43").Style = $__style
$ws.Range("E43").Value = '  -0.40%  '
$ws.Range("E44").Value = '  -5.33%  '
$__style = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '40.17'
$ws.Range("D46").Style = $__style
$ws.Range("E46").Value = '  -2.21%  '
$ws.Range("E47").Value = '  -4.19%  '
$__style = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.595'
$ws.Range("D48").Style = $__style
$ws.Range("E48").Value = '  -2.01%  '
$__style = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '153.33'
$ws.Range("D49").Style = $__style
$ws.Range("E49").Value = '  -1.97%  '
$ws.Range("E50").Value = '  -2.58%  '
$__style = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.72'
$ws.Range("D51").Style = $__style
$ws.Range("E51").Value = '  -4.36%  '
